$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# xlWhole = 1 -> only replace cells whose entire content matches
$xlWhole = 1

# Correct status_label: "bleu" -> "noir"
$used.Replace("bleu", "noir", $xlWhole)

# Correct status_name wording ("posté" -> "postés"/"publiés" phrasing).
# Replace the longer (36 mois) variant first so the shorter phrase below
# doesn't also need to re-match it.
$used.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $xlWhole)
$used.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $xlWhole)
$used.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $xlWhole)
